$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the seventh lesson: new row 9 with date 43046 and 2 hours.
# Copy the date cell's style/number format from the row above (A8) so it
# renders as a date just like the other entries.
$ws.Range("A8").Copy() | Out-Null
$ws.Range("A9").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A9").Value = 43046
$ws.Range("B9").Value = 2

# Update the sheet's dimension/used range naturally follows from the new data.

# Update the selection to match the new active cell reported in the diff.
$ws.Range("C9").Select() | Out-Null

$wb.Application.Calculate() | Out-Null
